# Updated For New Mail Sequences
# Adds three new Office365 asset entries to the "Assets" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$newAssets = @(
    "Office365ApplicationID",
    "Office365ApplicationSecret",
    "Office365TenantID"
)

$row = 17
foreach ($name in $newAssets) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $row = $row + 1
}

[void]$ws.Range("C17").Select()
